$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 275
$ws1.Range("F7").Value = 55
$ws1.Range("F10").Value = 101
$ws1.Range("F11").Value = 4332
$ws1.Range("C15").Value = " 江西·JMG（江西广电）第二届UP动漫游戏博览会-火只木南专场见面会"
$ws1.Range("F16").Value = 109
$ws1.Range("F20").Value = 3097
$ws1.Range("G20").Value = 59.9
$ws1.Range("F26").Value = 79
$ws1.Range("F30").Value = 199
$ws1.Range("F32").Value = 499
$ws1.Range("F33").Value = 1713
$ws1.Range("F34").Value = 260

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 275
$ws4.Range("F7").Value = 55
$ws4.Range("F10").Value = 101
$ws4.Range("F11").Value = 4332
$ws4.Range("C15").Value = " 江西·JMG（江西广电）第二届UP动漫游戏博览会-火只木南专场见面会"
$ws4.Range("F16").Value = 109
$ws4.Range("F20").Value = 3097
$ws4.Range("G20").Value = 59.9
$ws4.Range("F26").Value = 79
$ws4.Range("F30").Value = 199
$ws4.Range("F32").Value = 499
$ws4.Range("F33").Value = 1713
$ws4.Range("F34").Value = 260
